$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C38").Value2 = '[name="W"]  So what’s the deal with ''the boss of this place?''' + "`n"
$ws.Range("C116").Value2 = '[name="Theresa"]  Maybe it’s just my whimsy. The Doctor and Kal''tsit may reject its ''true name.''' + "`n"
$ws.Range("C120").Value2 = '[name="Theresa"]  ''Rhodes Island.''' + "`n"
$ws.Range("C140").Value2 = '[name="W"]  ''W...''' + "`n"
$ws.Range("C141").Value2 = '[name="Theresa"]  That’s not what I mean. ''W'' is a mercenary codename. I want to know your true name. The one that represents you..' + "`n"
$ws.Range("C149").Value2 = '[name="Theresa"]  If you wait for the dust to settle in Kazdel... until you’re no longer ''W,'' we could speak like this once again.' + "`n"
$ws.Range("C178").Value2 = '[name="W"]  (That hood... is that the ''Doctor'' they were talking about?)' + "`n"
$ws.Range("C184").Value2 = '[name="W"]  (''The Doctor,'' huh?)' + "`n"
